# Add a new "Player Info" worksheet as the first sheet in the workbook,
# populate it with player metadata, and update the MATCH_CARD_LINK
# columns on the existing sheets to hold just the numeric MATCH_CODE.

$wb = $excel.ActiveWorkbook

# --- Add the new "Player Info" sheet before the existing first sheet ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$newSheet = $wb.Worksheets.Add($battingSheet)
$newSheet.Name = "Player Info"

# NOTE: references obtained before the Add() call above can become stale
# (the collection can resolve them positionally), so re-fetch every sheet
# we still need to touch by name now that the sheet collection is final.
$playerInfo = $wb.Worksheets.Item("Player Info")
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Match the page margins used by the rest of the workbook (0.75in / 1in /
# 0.5in) instead of the engine's default template for brand-new sheets.
$playerInfo.PageSetup.LeftMargin = 54
$playerInfo.PageSetup.RightMargin = 54
$playerInfo.PageSetup.TopMargin = 72
$playerInfo.PageSetup.BottomMargin = 72
$playerInfo.PageSetup.HeaderMargin = 36
$playerInfo.PageSetup.FooterMargin = 36

# Headers (bold / centered / bordered style, matching the other header rows).
# Copy/PasteSpecial(Formats) reuses the existing header cell style instead
# of fabricating a new (duplicate) style entry.
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$battingSheet.Range("A1:D1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Data row. "ID" is purely numeric-looking text, so force a Text number
# format first - otherwise Excel auto-coerces the assigned string into a
# real number, which would change the stored cell type.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5949"
$playerInfo.Range("B2").Value = "Riley Patrick Meredith"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast"

# --- Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4484"

# --- Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4484"
